# Adicionando as informações sobre idade materna
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns O, P, Q (row 1)
$ws.Range("O1").Value = "Média de idade materna"
$ws.Range("P1").Value = "Idade materna com  <= 15 anos"
$ws.Range("Q1").Value = "Idade materna com  >= 35 anos"

# Match the formatting already used by the rest of the header row (bold + centered)
$headerRange = $ws.Range("O1:Q1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter

# New data values per year row (2012-2023 => rows 2-13)
$data = @{
    2  = @(26.91, 10917, 85168)
    3  = @(27.01, 11039, 87050)
    4  = @(27.16, 11209, 94104)
    5  = @(27.34, 10343, 99964)
    6  = @(27.49, 8859,  99693)
    7  = @(27.77, 7909,  109014)
    8  = @(28.03, 7093,  115692)
    9  = @(28.19, 6207,  115337)
    10 = @(28.27, 5574,  110781)
    11 = @(28.37, 5138,  107904)
    12 = @(28.59, 4189,  109063)
    13 = @(28.62, 3944,  107376)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 15).Value = $vals[0]  # Coluna O - Média de idade materna
    $ws.Cells.Item($row, 16).Value = $vals[1]  # Coluna P - Idade materna <= 15 anos
    $ws.Cells.Item($row, 17).Value = $vals[2]  # Coluna Q - Idade materna >= 35 anos
}
